$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (30) down into the
# new row (31), matching columns A:K only (row 30's L column is intentionally
# left out of the new row).
$ws.Range("A30:K30").Copy() | Out-Null
$ws.Range("A31:K31").PasteSpecial(-4122) | Out-Null

# Populate the new admin code entry ("two mid digits" case/code variant).
$ws.Range("A31").Value = "22TRC00571"
$ws.Range("B31").Value = "22TRC00571-A"
$ws.Range("C31").Value = "Mick"
$ws.Range("D31").Value = "Jagger"
$ws.Range("E31").Value = "Test TWO"
$ws.Range("F31").Value = "1501.17-11-01"
$ws.Range("G31").Value = "MM"
$ws.Range("H31").Value = "Y"
$ws.Range("I31").Value = $true

# Move the active selection to the newly added row, as in the source workbook.
$ws.Range("G31").Select() | Out-Null
